$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.481.44"
$ws.Range("E2").Value = "  +3.78%  "
$ws.Range("D3").Value = "3.085.38"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("D4").Value = "'0.996"
$ws.Range("D4").Style = $ws.Range("D4").Style
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'218.78"
$ws.Range("D5").Style = $ws.Range("D5").Style
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "'618.02"
$ws.Range("D6").Style = $ws.Range("D6").Style
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("D7").Value = "'0.376"
$ws.Range("D7").Style = $ws.Range("D7").Style
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("D8").Value = "'0.922"
$ws.Range("D8").Style = $ws.Range("D8").Style
$ws.Range("E8").Value = "  +12.56%  "
$ws.Range("D9").Value = "'0.998"
$ws.Range("D9").Style = $ws.Range("D9").Style
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "3.090.93"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "'0.679"
$ws.Range("D11").Style = $ws.Range("D11").Style
$ws.Range("E11").Value = "  +16.92%  "
$ws.Range("E12").Value = "  +6.18%  "
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("D13").Style = $ws.Range("D13").Style
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "91.196.70"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "'33.01"
$ws.Range("D16").Style = $ws.Range("D16").Style
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "3.657.23"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "3.080.58"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'3.48"
$ws.Range("D19").Style = $ws.Range("D19").Style
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "'0.0000222"
$ws.Range("D20").Style = $ws.Range("D20").Style
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  +4.04%  "
$ws.Range("D22").Value = "'435.28"
$ws.Range("D22").Style = $ws.Range("D22").Style
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "'8.47"
$ws.Range("D23").Style = $ws.Range("D23").Style
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'5.16"
$ws.Range("D24").Style = $ws.Range("D24").Style
$ws.Range("E24").Value = "  +5.35%  "
$ws.Range("D25").Value = "'5.62"
$ws.Range("D25").Style = $ws.Range("D25").Style
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").Value = "'84.19"
$ws.Range("D26").Style = $ws.Range("D26").Style
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").Value = "'11.81"
$ws.Range("D27").Style = $ws.Range("D27").Style
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").Value = "3.257.48"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +7.19%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "'8.83"
$ws.Range("D32").Style = $ws.Range("D32").Style
$ws.Range("E32").Value = "  +8.31%  "
$ws.Range("D33").Value = "'3.90"
$ws.Range("D33").Style = $ws.Range("D33").Style
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "'520.82"
$ws.Range("D34").Style = $ws.Range("D34").Style
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("D35").Value = "'7.11"
$ws.Range("D35").Style = $ws.Range("D35").Style
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'1.28"
$ws.Range("D36").Style = $ws.Range("D36").Style
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.138"
$ws.Range("D37").Style = $ws.Range("D37").Style
$ws.Range("E37").Value = "  -7.13%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'23.02"
$ws.Range("D39").Style = $ws.Range("D39").Style
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = $ws.Range("D41").Style
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'0.142"
$ws.Range("D43").Style = $ws.Range("D43").Style
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.88"
$ws.Range("D44").Style = $ws.Range("D44").Style
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").Value = "'0.368"
$ws.Range("D45").Style = $ws.Range("D45").Style
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "'0.0724"
$ws.Range("D46").Style = $ws.Range("D46").Style
$ws.Range("E46").Value = "  +10.23%  "
$ws.Range("D47").Value = "'43.86"
$ws.Range("D47").Style = $ws.Range("D47").Style
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'141.11"
$ws.Range("D48").Style = $ws.Range("D48").Style
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("E49").Value = "  +13.22%  "
$ws.Range("E50").Value = "  +6.34%  "
$ws.Range("D51").Value = "'164.67"
$ws.Range("D51").Style = $ws.Range("D51").Style
$ws.Range("E51").Value = "  +1.24%  "
